$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1125.909
$ws.Range("I98").Value = 1134
$ws.Range("K98").Value = 1134
$ws.Range("M98").Value = 364
$ws.Range("H112").Value = 3518.3125
$ws.Range("I112").Value = 949
$ws.Range("J112").Value = 3885.3572
$ws.Range("K112").Value = 2847
$ws.Range("L112").Value = 11656.0716
$ws.Range("M112").Value = -1739
$ws.Range("N112").Value = -13872.0716
$ws.Range("H113").Value = 3299.8
$ws.Range("I113").Value = 3296
$ws.Range("J113").Value = 3305.5
$ws.Range("K113").Value = 3296
$ws.Range("L113").Value = 3305.5
$ws.Range("M113").Value = -42
$ws.Range("N113").Value = -9813.5
$ws.Range("H116").Value = 5796
$ws.Range("I116").Value = 4180.8
$ws.Range("J116").Value = 6221.0527
$ws.Range("K116").Value = 4180.8
$ws.Range("L116").Value = 6221.0527
$ws.Range("M116").Value = -738.8000000000002
$ws.Range("N116").Value = -13105.0527
$ws.Range("H122").Value = 1125.909
$ws.Range("I122").Value = 1134
$ws.Range("K122").Value = 3402
$ws.Range("M122").Value = -952
$ws.Range("H125").Value = 5267.5
$ws.Range("J125").Value = 9036
$ws.Range("L125").Value = 81324
$ws.Range("N125").Value = -86244
$ws.Range("H132").Value = 4628.96
$ws.Range("I132").Value = 1419.6666
$ws.Range("J132").Value = 7591.385
$ws.Range("K132").Value = 4258.9998
$ws.Range("L132").Value = 22774.155
$ws.Range("M132").Value = -1728.9998
$ws.Range("N132").Value = -27834.155
$ws.Range("H135").Value = 3611.1667
$ws.Range("I135").Value = 3510.5
$ws.Range("K135").Value = 31594.5
$ws.Range("M135").Value = -29059.5
$ws.Range("H141").Value = 10799.2
$ws.Range("I141").Value = 10799.2
$ws.Range("K141").Value = 32397.6
$ws.Range("M141").Value = -27217.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4729.7856
$ws.Range("I32").Value = 4081.963
$ws.Range("J32").Value = 22221
$ws.Range("K32").Value = 4081.963
$ws.Range("L32").Value = 22221
$ws.Range("M32").Value = -3794.963
$ws.Range("N32").Value = -22795
$ws.Range("H110").Value = 2200
$ws.Range("I110").Value = 2200
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 2200
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -155
$ws.Range("N110").ClearContents()
$ws.Range("H132").Value = 4260
$ws.Range("I132").Value = 4080
$ws.Range("K132").Value = 12240
$ws.Range("M132").Value = -9710

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3538.5
$ws.Range("I86").Value = 3831.7778
$ws.Range("K86").Value = 3831.7778
$ws.Range("M86").Value = -2708.7778
$ws.Range("H89").Value = 3538.5
$ws.Range("I89").Value = 3831.7778
$ws.Range("K89").Value = 19158.889
$ws.Range("M89").Value = -13542.889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1595.8462
$ws.Range("I7").Value = 1467
$ws.Range("J7").Value = 1746.1666
$ws.Range("K7").Value = 1467
$ws.Range("L7").Value = 1746.1666
$ws.Range("M7").Value = -1354
$ws.Range("N7").Value = -1972.1666
$ws.Range("H16").Value = 824.75
$ws.Range("I16").Value = 771.1429000000001
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 771.1429000000001
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -484.1429000000001
$ws.Range("N16").Value = -1774
$ws.Range("H31").Value = 4873.6665
$ws.Range("I31").Value = 3692.7144
$ws.Range("K31").Value = 3692.7144
$ws.Range("M31").Value = -3397.7144
$ws.Range("H34").Value = 4873.6665
$ws.Range("I34").Value = 3692.7144
$ws.Range("K34").Value = 3692.7144
$ws.Range("M34").Value = -3490.7144
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()
$ws.Range("H113").Value = 824.75
$ws.Range("I113").Value = 771.1429000000001
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 771.1429000000001
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1398.8571
$ws.Range("N113").Value = -5540
$ws.Range("H132").Value = 7806.625
$ws.Range("I132").Value = 6146.654
$ws.Range("K132").Value = 18439.962
$ws.Range("M132").Value = -15909.962
$ws.Range("H134").Value = 2241.8096
$ws.Range("I134").Value = 2112.6667
$ws.Range("J134").Value = 3016.6667
$ws.Range("K134").Value = 6338.000100000001
$ws.Range("L134").Value = 9050.000100000001
$ws.Range("M134").Value = -3803.000100000001
$ws.Range("N134").Value = -14120.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 904.1429000000001
$ws.Range("I26").Value = 796.75
$ws.Range("J26").Value = 1047.3334
$ws.Range("K26").Value = 2390.25
$ws.Range("L26").Value = 3142.0002
$ws.Range("M26").Value = -2102.25
$ws.Range("N26").Value = -3718.0002
$ws.Range("H92").Value = 874.1667
$ws.Range("J92").Value = 1000
$ws.Range("L92").Value = 3000
$ws.Range("N92").Value = -5496
$ws.Range("H131").Value = 993.5
$ws.Range("I131").Value = 993.5
$ws.Range("K131").Value = 2980.5
$ws.Range("M131").Value = 2059.5
$ws.Range("H140").Value = 853.6667
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 102.14286
$ws.Range("I2").Value = 8.916667
$ws.Range("J2").Value = 226.44444
$ws.Range("K2").Value = 8.916667
$ws.Range("L2").Value = 226.44444
$ws.Range("M2").Value = 104.083333
$ws.Range("N2").Value = -452.44444
$ws.Range("H107").Value = 694.1111
$ws.Range("I107").Value = 530.875
$ws.Range("K107").Value = 530.875
$ws.Range("M107").Value = 1389.125
$ws.Range("H113").Value = 2405.5
$ws.Range("I113").Value = 2405.5
$ws.Range("K113").Value = 2405.5
$ws.Range("M113").Value = -235.5
$ws.Range("H122").Value = 3142.5
$ws.Range("I122").Value = 3155.4546
$ws.Range("K122").Value = 9466.363799999999
$ws.Range("M122").Value = -7016.363799999999
$ws.Range("H132").Value = 1941.7142
$ws.Range("I132").Value = 1879.7106
$ws.Range("J132").Value = 2530.75
$ws.Range("K132").Value = 5639.1318
$ws.Range("L132").Value = 7592.25
$ws.Range("M132").Value = -3109.1318
$ws.Range("N132").Value = -12652.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1704.7778
$ws.Range("J22").Value = 1894.25
$ws.Range("L22").Value = 1894.25
$ws.Range("N22").Value = -2484.25
$ws.Range("H27").Value = 1704.7778
$ws.Range("J27").Value = 1894.25
$ws.Range("L27").Value = 1894.25
$ws.Range("N27").Value = -2108.25
$ws.Range("H76").Value = 14000
$ws.Range("J76").Value = 14000
$ws.Range("L76").Value = 14000
$ws.Range("N76").Value = -14676
$ws.Range("H79").Value = 14000
$ws.Range("J79").Value = 14000
$ws.Range("L79").Value = 14000
$ws.Range("N79").Value = -16340
$ws.Range("H132").Value = 2743.5557
$ws.Range("I132").Value = 2711.625
$ws.Range("K132").Value = 8134.875
$ws.Range("M132").Value = -5604.875
$ws.Range("H136").Value = 3782.1667
$ws.Range("I136").Value = 3665.3333
$ws.Range("K136").Value = 10995.9999
$ws.Range("M136").Value = -8445.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 6972877
$ws.Range("I100").Value = 34848984
$ws.Range("J100").Value = 3850
$ws.Range("K100").Value = 69697968
$ws.Range("L100").Value = 7700
$ws.Range("M100").Value = -69697427
$ws.Range("N100").Value = -8782
$ws.Range("H132").Value = 112519.555
$ws.Range("I132").Value = 126397
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 379191
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -376661
$ws.Range("N132").Value = -9560
$ws.Range("H136").Value = 2192.9412
$ws.Range("I136").Value = 1505.7142
$ws.Range("K136").Value = 4517.142599999999
$ws.Range("M136").Value = -1967.142599999999
